# Applies the "finish from my side" edit to slide 2 (File Structure Figures):
#   - Shape id 45 (".csv files of loss documentation", top-left cell) is
#     widened/shifted left and its text changed from "... of ..." to "... for ...".
#   - The eight other ".csv files of loss documentation" cells (ids 47, 49,
#     54, 56, 58, 63, 65, 67) are removed.
#   - Ten replacement cells (same style/text, nudged positions) are added
#     back so every Cycle-GAN block again has a ".csv files for loss
#     documentation" cell.

function EmuToPt($emu) {
    # PowerPoint COM measures Left/Top/Width/Height in points (1 pt = 12700 EMU).
    # The host's pt -> EMU conversion truncates rather than rounds, so nudge
    # the point value up by half an EMU to land exactly on the target EMU.
    return ($emu + 0.5) / 12700
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# 1) Modify shape 45: reposition/resize and fix its text.
$sh45 = Get-ShapeById $s 45
$sh45.Left = EmuToPt 335975
$sh45.Width = EmuToPt 1591540
$sh45.TextFrame.TextRange.Text = ".csv files for loss documentation"

# 2) Delete the other eight "csv files of loss documentation" cells.
$deleteIds = @(47, 49, 54, 56, 58, 63, 65, 67)
foreach ($id in $deleteIds) {
    $sh = Get-ShapeById $s $id
    if ($sh -ne $null) {
        $sh.Delete()
    }
}

# 3) Re-add ten replacement cells by duplicating the fixed-up shape 45 and
#    moving each copy into place.
$newShapes = @(
    @(335975, 3408218, "Rechteck 32"),
    @(335975, 5243945, "Rechteck 33"),
    @(335975, 1544783, "Rechteck 34"),
    @(335975, 3408219, "Rechteck 35"),
    @(4379338, 5240484, "Rechteck 36"),
    @(4379338, 1541322, "Rechteck 37"),
    @(4379338, 3404758, "Rechteck 38"),
    @(8407142, 5240484, "Rechteck 68"),
    @(8407142, 1541322, "Rechteck 69"),
    @(8407142, 3404758, "Rechteck 70")
)

foreach ($entry in $newShapes) {
    $dup = $sh45.Duplicate()
    $dup.Left = EmuToPt $entry[0]
    $dup.Top = EmuToPt $entry[1]
    $dup.Width = EmuToPt 1591540
    $dup.Height = EmuToPt 1184563
    $dup.Name = $entry[2]
}
